{"js": "// Replace the 25 two-digit-division problems in the single table with the\n// new set of problems, preserving cell/run formatting. The problems occupy\n// the table's 1st, 5th, 9th, 13th, and 17th rows (0-based: 0, 4, 8, 12, 16);\n// the rows in between are blank spacer rows left untouched.\nconst newProblems = [\n  \"98\u00f78=\", \"56\u00f78=\", \"30\u00f73=\", \"44\u00f72=\", \"37\u00f77=\",\n  \"31\u00f76=\", \"97\u00f73=\", \"23\u00f73=\", \"87\u00f75=\", \"32\u00f73=\",\n  \"43\u00f75=\", \"93\u00f73=\", \"95\u00f72=\", \"43\u00f78=\", \"12\u00f78=\",\n  \"55\u00f77=\", \"51\u00f75=\", \"86\u00f74=\", \"77\u00f75=\", \"77\u00f78=\",\n  \"77\u00f77=\", \"86\u00f75=\", \"84\u00f76=\", \"21\u00f74=\", \"64\u00f79=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst problemRowIndexes = [0, 4, 8, 12, 16];\nlet k = 0;\nfor (const rowIdx of problemRowIndexes) {\n  for (let col = 0; col < 5; col++) {\n    table.getCell(rowIdx, col).value = newProblems[k];\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit-division problems in the single table with the\n# new set of problems, preserving cell/run formatting. The problems occupy\n# the table's rows 1, 5, 9, 13, and 17 (1-based); the rows in between are\n# blank spacer rows left untouched.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = @(1, 5, 9, 13, 17)\n$values = @(\n    @(\"98\u00f78=\", \"56\u00f78=\", \"30\u00f73=\", \"44\u00f72=\", \"37\u00f77=\"),\n    @(\"31\u00f76=\", \"97\u00f73=\", \"23\u00f73=\", \"87\u00f75=\", \"32\u00f73=\"),\n    @(\"43\u00f75=\", \"93\u00f73=\", \"95\u00f72=\", \"43\u00f78=\", \"12\u00f78=\"),\n    @(\"55\u00f77=\", \"51\u00f75=\", \"86\u00f74=\", \"77\u00f75=\", \"77\u00f78=\"),\n    @(\"77\u00f77=\", \"86\u00f75=\", \"84\u00f76=\", \"21\u00f74=\", \"64\u00f79=\")\n)\n\nfor ($i = 0; $i -lt $rows.Count; $i++) {\n    $r = $rows[$i]\n    $rowVals = $values[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowVals[$c - 1]\n    }\n}\n"}
